$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.883.18'
$ws.Range("E2").Value = '  +0.63%  '

$ws.Range("D3").Value = '1.636.11'
$ws.Range("E3").Value = '  +0.16%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.37%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.00'
$ws.Range("E5").Value = '  -0.12%  '

$ws.Range("E6").Value = '  -0.24%  '

$ws.Range("E7").Value = '  -0.34%  '

$ws.Range("E8").Value = '  -0.73%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0632'
$ws.Range("E9").Value = '  -0.59%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.68'
$ws.Range("E10").Value = '  +0.42%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0791'
$ws.Range("E11").Value = '  +0.67%  '

$ws.Range("E12").Value = '  +0.56%  '

$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.659.05'
$ws.Range("E13").Value = '  +1.51%  '

$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '1.862.19'
$ws.Range("E14").Value = '  +0.21%  '

$ws.Range("E15").Value = '  -0.67%  '

$ws.Range("D16").Value = '0.0₃0759'
$ws.Range("E16").Value = '  -0.70%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.91'
$ws.Range("E17").Value = '  +0.36%  '

$ws.Range("D18").Value = '25.884.21'
$ws.Range("E18").Value = '  +0.50%  '

$ws.Range("E19").Value = '  -0.25%  '

$ws.Range("E20").Value = '  +0.09%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '191.66'
$ws.Range("E21").Value = '  -1.10%  '

$ws.Range("E22").Value = '  +0.43%  '

$ws.Range("E23").Value = '  +0.78%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.00'
$ws.Range("E24").Value = '  -0.28%  '

$ws.Range("E25").Value = '  -1.80%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '142.47'
$ws.Range("E26").Value = '  +0.48%  '

$ws.Range("E27").Value = '  +1.18%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.86'
$ws.Range("E28").Value = '  -0.31%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.52'
$ws.Range("E29").Value = '  +0.05%  '

$ws.Range("E30").Value = '  -0.43%  '

$ws.Range("E31").Value = '  +0.63%  '

$ws.Range("E32").Value = '  +0.34%  '

$ws.Range("E33").Value = '  +0.00%  '

$ws.Range("E34").Value = '  +0.95%  '

$ws.Range("E35").Value = '  +0.43%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.908'
$ws.Range("E36").Value = '  +0.86%  '

$ws.Range("D37").Value = '1.148.31'
$ws.Range("E37").Value = '  +2.23%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.545'
$ws.Range("E38").Value = '  -0.52%  '

$ws.Range("E39").Value = '  -0.97%  '

$ws.Range("E40").Value = '  +0.65%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'

$ws.Range("E42").Value = '  +1.14%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '100.57'
$ws.Range("E43").Value = '  +0.89%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.803'
$ws.Range("E44").Value = '  -0.08%  '

$ws.Range("D45").Value = '1.771.94'
$ws.Range("E45").Value = '  +0.21%  '

$ws.Range("E46").Value = '  -1.49%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '55.64'
$ws.Range("E47").Value = '  +1.04%  '

$ws.Range("E48").Value = '  +1.97%  '

$ws.Range("E49").Value = '  +5.78%  '

$ws.Range("E50").Value = '  -0.03%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.60'
$ws.Range("E51").Value = '  +0.55%  '
